$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 82
$ws.Range("E2").Value = 0.803921568627451
$ws.Range("F2").Value = 0.803921568627451
$ws.Range("G2").Value = 0.1020570907636317
$ws.Range("H2").Value = 0.08204589649625291
$ws.Range("I2").Value = 452494.6720494052
$ws.Range("J2").Value = 163633.3362257026
$ws.Range("L2").Value = 163633.3362257026
$ws.Range("M2").Value = 616128.0082751078
$ws.Range("N2").Value = 10292112.5688
$ws.Range("O2").Value = 9884371.638699999
$ws.Range("P2").Value = 0.01589890657839757
$ws.Range("Q2").Value = 0.01655475352474948

$ws.Range("C3").Value = 102
$ws.Range("D3").Value = 87
$ws.Range("E3").Value = 0.8529411764705882
$ws.Range("F3").Value = 0.8446601941747572
$ws.Range("G3").Value = 0.09750327982997359
$ws.Range("H3").Value = 0.08235713927386119
$ws.Range("I3").Value = 473488.5336105639
$ws.Range("J3").Value = 171531.8470123119
$ws.Range("L3").Value = 171531.8470123119
$ws.Range("M3").Value = 645020.3806228759
$ws.Range("N3").Value = 10436953.712764
$ws.Range("O3").Value = 10029580.554761
$ws.Range("P3").Value = 0.01643504912765254
$ws.Range("Q3").Value = 0.01710259427856995

$ws.Range("D4").Value = 87
$ws.Range("E4").Value = 0.8365384615384616
$ws.Range("F4").Value = 0.8365384615384616
$ws.Range("G4").Value = 0.09749531305237304
$ws.Range("H4").Value = 0.08155857918804287
$ws.Range("I4").Value = 497757.4148347
$ws.Range("J4").Value = 176579.9150305909
$ws.Range("L4").Value = 176579.9150305909
$ws.Range("M4").Value = 674337.3298652909
$ws.Range("N4").Value = 10911311.62264692
$ws.Range("O4").Value = 10502967.26990383
$ws.Range("P4").Value = 0.01618319787183892
$ws.Range("Q4").Value = 0.01681238363339275

$ws.Range("D5").Value = 86
$ws.Range("E5").Value = 0.819047619047619
$ws.Range("F5").Value = 0.819047619047619
$ws.Range("G5").Value = 0.09816608625862366
$ws.Range("H5").Value = 0.08040269922134892
$ws.Range("I5").Value = 513130.3226945847
$ws.Range("J5").Value = 181692.4051889305
$ws.Range("L5").Value = 181692.4051889305
$ws.Range("M5").Value = 694822.7278835152
$ws.Range("N5").Value = 11330655.55062633
$ws.Range("O5").Value = 10919960.86730095
$ws.Range("P5").Value = 0.01603547159095195
$ws.Range("Q5").Value = 0.01663855826928791

$ws.Range("D6").Value = 88
$ws.Range("E6").Value = 0.8301886792452831
$ws.Range("F6").Value = 0.8301886792452831
$ws.Range("G6").Value = 0.09702974602399582
$ws.Range("H6").Value = 0.08055299669916634
$ws.Range("I6").Value = 535915.2845663553
$ws.Range("J6").Value = 190023.7034400649
$ws.Range("L6").Value = 190023.7034400649
$ws.Range("M6").Value = 725938.9880064201
$ws.Range("N6").Value = 11688284.31584512
$ws.Range("O6").Value = 11273818.79201997
$ws.Range("P6").Value = 0.01625762158971963
$ws.Range("Q6").Value = 0.01685530936283726

